$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# "About" sheet
$aboutWs = $wb.Worksheets.Item("About")
$aboutWs.Range("A2").Value = "Version: " + $newVersion
$aboutWs.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Shuiyu Coal Mine, China, M1201, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

# "Boundaries and methane sources" sheet
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 8; $r++) {
    $dataWs.Cells.Item($r, 19).Value = $newVersion
}
